$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "Pediatrics" to "Session"
$ws.Name = "Session"

# Remove the second log entry (row 3) entirely, shifting rows up
$ws.Rows.Item(3).Delete()
